$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-06 Sunday", "2025-04-07 Monday"),
    @("766×4=3064", "635×7=4445"),
    @("811×7=5677", "258×9=2322"),
    @("801×9=7209", "289×5=1445"),
    @("345×9=3105", "272×8=2176"),
    @("572×2=1144", "826×4=3304"),
    @("904×7=6328", "116×4=464"),
    @("332×5=1660", "267×6=1602"),
    @("575×7=4025", "419×4=1676"),
    @("914×4=3656", "900×5=4500"),
    @("704×9=6336", "445×3=1335"),
    @("420×3=1260", "136×4=544"),
    @("866×4=3464", "826×7=5782"),
    @("624×5=3120", "862×4=3448"),
    @("251×2=502", "108×9=972"),
    @("352×3=1056", "572×8=4576"),
    @("122×3=366", "569×4=2276"),
    @("432×5=2160", "957×5=4785"),
    @("883×2=1766", "303×6=1818"),
    @("765×9=6885", "628×7=4396"),
    @("308×5=1540", "773×7=5411"),
    @("378×4=1512", "645×3=1935"),
    @("684×6=4104", "275×9=2475"),
    @("359×3=1077", "691×3=2073"),
    @("899×6=5394", "652×2=1304"),
    @("230×9=2070", "481×5=2405")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
